# Updated test cases which were failing
# Refresh the "Date" column (column B) timestamps on the rows whose
# tests were re-run, across the affected sheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "VerifySuccessfulPaymentCCNoCF";    Cell = "B2"; Value = "Wed Oct 15 19:45:21 IST 2025" },
    @{ Sheet = "VerifySuccessfulPaymentCCNoCF";    Cell = "B3"; Value = "Wed Oct 15 19:48:03 IST 2025" },
    @{ Sheet = "VerifySuccessfulPaymentCCNoCF";    Cell = "B4"; Value = "Wed Oct 15 19:50:37 IST 2025" },

    @{ Sheet = "CreateVerifyDeleteCCPM";           Cell = "B2"; Value = "Wed Oct 15 19:33:04 IST 2025" },

    @{ Sheet = "VerifyCFVerbiageOnRPCCDCF";        Cell = "B2"; Value = "Wed Oct 15 19:36:56 IST 2025" },

    @{ Sheet = "VerifyStaticTextOnRecieptCC";      Cell = "B2"; Value = "Wed Oct 15 19:41:23 IST 2025" },

    @{ Sheet = "VerifySuccessfulPaymentCCSCF";     Cell = "B2"; Value = "Wed Oct 15 19:53:23 IST 2025" },
    @{ Sheet = "VerifySuccessfulPaymentCCSCF";     Cell = "B3"; Value = "Wed Oct 15 19:56:17 IST 2025" },
    @{ Sheet = "VerifySuccessfulPaymentCCSCF";     Cell = "B4"; Value = "Wed Oct 15 19:59:00 IST 2025" },

    @{ Sheet = "VerifySuccessfulPaymentCCDCF";     Cell = "B2"; Value = "Wed Oct 15 19:43:48 IST 2025" },
    @{ Sheet = "VerifySuccessfulPaymentCCDCF";     Cell = "B3"; Value = "Wed Oct 15 19:44:17 IST 2025" },
    @{ Sheet = "VerifySuccessfulPaymentCCDCF";     Cell = "B4"; Value = "Wed Oct 15 19:44:46 IST 2025" },

    @{ Sheet = "VerifySuccessfulPaymentSPMCCDCF";  Cell = "B2"; Value = "Wed Oct 15 00:42:59 IST 2025" },

    @{ Sheet = "UiVerificationSPPaymentMethodCC";  Cell = "B2"; Value = "Tue Oct 14 17:32:15 IST 2025" },

    @{ Sheet = "VerifyStaticTextOnPPCCDCF";        Cell = "B2"; Value = "Wed Oct 15 19:39:06 IST 2025" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
